$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host ("P71(sum) format: " + $ws.Cells.Item(71,16).NumberFormat)
Write-Host ("P71(sum) value: " + $ws.Cells.Item(71,16).Value)
Write-Host ("A72(footer) format: " + $ws.Cells.Item(72,1).NumberFormat)
Write-Host ("A72(footer) value: " + $ws.Cells.Item(72,1).Value)
Write-Host ("G72(footer) value: " + $ws.Cells.Item(72,7).Value)
Write-Host ("K72(footer) value: " + $ws.Cells.Item(72,11).Value)
